$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Final values for B2:J49 (columns B..J, rows 2..49)
$data = @(
    @(86,46,44,3,23,39,58,16,29),
    @(6,4,3,0,1,1,7,0,4),
    @(5,4,0,0,0,0,3,0,0),
    @(0,0,3,0,1,0,3,0,3),
    @(1,0,0,0,0,1,1,0,1),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(10,8,7,0,4,4,14,2,6),
    @(4,4,0,0,3,2,3,0,0),
    @(0,0,2,0,1,1,1,1,1),
    @(2,1,2,0,0,1,1,0,1),
    @(0,0,3,0,0,0,6,0,4),
    @(4,3,0,0,0,0,3,1,0),
    @(11,5,7,3,6,8,11,1,5),
    @(3,3,4,0,1,3,1,1,1),
    @(0,0,1,0,0,0,4,0,2),
    @(6,2,1,3,5,4,3,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(2,0,0,0,0,1,1,0,0),
    @(0,0,1,0,0,0,1,0,2),
    @(31,11,17,0,6,9,6,4,2),
    @(4,2,4,0,1,3,1,0,2),
    @(4,1,2,0,0,0,3,0,0),
    @(11,3,2,0,3,0,0,1,0),
    @(8,2,2,0,0,3,0,0,0),
    @(2,0,2,0,0,0,0,0,0),
    @(2,3,5,0,2,3,2,3,0),
    @(4,4,4,0,2,5,5,3,4),
    @(0,0,0,0,0,0,0,0,0),
    @(1,1,1,0,0,1,1,0,1),
    @(1,0,2,0,2,4,2,3,2),
    @(0,0,0,0,0,0,0,0,0),
    @(2,3,0,0,0,0,2,0,1),
    @(0,0,1,0,0,0,0,0,0),
    @(13,8,3,0,2,8,8,3,4),
    @(0,0,0,0,0,0,0,0,0),
    @(4,3,0,0,2,0,2,1,3),
    @(3,3,2,0,0,2,0,1,1),
    @(3,0,0,0,0,0,0,1,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,4,0,0),
    @(3,2,1,0,0,6,2,0,0),
    @(11,6,3,0,2,4,7,3,4),
    @(2,1,1,0,0,1,0,0,0),
    @(3,2,0,0,0,0,0,1,0),
    @(2,0,1,0,0,2,6,0,2),
    @(4,3,1,0,1,0,1,1,0),
    @(0,0,0,0,1,1,0,1,2)
)

$startRow = 2
$startCol = 2  # column B

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $c = $startCol + $j
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}
